$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A24").Value = "25/08/2023"
$ws.Range("B24").Value = "Add User form"
$ws.Range("C24").Value = "grid data is not shown on controls esp textboxes"

$ws.Range("C24").Select()
